# Refresh the KHL "Probabilities" tour sheet: rows 2-4 get the next
# batch of matches (replacing the 2025-11-30 fixtures) and a new row 5
# is appended for the additional 2025-12-01 fixture.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2
$ws.Cells.Item(2, 1).Value = 1369
$ws.Cells.Item(2, 2).Value = "2025-12-01T15:30:00"
$ws.Cells.Item(2, 3).Value = "Сибирь"
$ws.Cells.Item(2, 4).Value = "Барыс"
$ws.Cells.Item(2, 5).Value = 897836
$ws.Cells.Item(2, 6).Value = "https://text.khl.ru/text/897836.html"
$ws.Cells.Item(2, 7).Value = 0.846154
$ws.Cells.Item(2, 8).Value = 1.854302
$ws.Cells.Item(2, 9).Value = 3.846402
$ws.Cells.Item(2, 10).Value = 4.961538
$ws.Cells.Item(2, 11).Value = 2.903846
$ws.Cells.Item(2, 12).Value = 2.850352
$ws.Cells.Item(2, 13).Value = 2.700456
$ws.Cells.Item(2, 14).Value = 23.063363
$ws.Cells.Item(2, 15).Value = 29.982381
$ws.Cells.Item(2, 16).Value = 53.045744
$ws.Cells.Item(2, 17).Value = -0.2
$ws.Cells.Item(2, 18).Value = -0.05563
$ws.Cells.Item(2, 19).Value = 0.423308
$ws.Cells.Item(2, 20).Value = 0.17033
$ws.Cells.Item(2, 21).Value = 0.405949
$ws.Cells.Item(2, 22).Value = 0.174522
$ws.Cells.Item(2, 23).Value = 0.825065
$ws.Cells.Item(2, 24).Value = 0.319303
$ws.Cells.Item(2, 25).Value = 0.680284
$ws.Cells.Item(2, 26).Value = 0.485923
$ws.Cells.Item(2, 27).Value = 0.513664
$ws.Cells.Item(2, 28).Value = 0.645717
$ws.Cells.Item(2, 29).Value = 0.35387
$ws.Cells.Item(2, 30).Value = 0.777072
$ws.Cells.Item(2, 31).Value = 0.222515
$ws.Cells.Item(2, 32).Value = 0.786022
$ws.Cells.Item(2, 33).Value = 0.213978
$ws.Cells.Item(2, 34).Value = 0.554926
$ws.Cells.Item(2, 35).Value = 0.445074
$ws.Cells.Item(2, 36).Value = 0.777357
$ws.Cells.Item(2, 37).Value = 0.222643
$ws.Cells.Item(2, 38).Value = 0.542462
$ws.Cells.Item(2, 39).Value = 0.457538
$ws.Cells.Item(2, 40).Value = 0.746927
$ws.Cells.Item(2, 41).Value = 0.732445

# Row 3
$ws.Cells.Item(3, 1).Value = 1369
$ws.Cells.Item(3, 2).Value = "2025-12-01T19:00:00"
$ws.Cells.Item(3, 3).Value = "Северсталь"
$ws.Cells.Item(3, 4).Value = "Трактор"
$ws.Cells.Item(3, 5).Value = 897835
$ws.Cells.Item(3, 6).Value = "https://text.khl.ru/text/897835.html"
$ws.Cells.Item(3, 7).Value = 1.46875
$ws.Cells.Item(3, 8).Value = 3.9375
$ws.Cells.Item(3, 9).Value = 1.15625
$ws.Cells.Item(3, 10).Value = 3.831021
$ws.Cells.Item(3, 11).Value = 2.649885
$ws.Cells.Item(3, 12).Value = 2.546875
$ws.Cells.Item(3, 13).Value = 5.40625
$ws.Cells.Item(3, 14).Value = 23.963798
$ws.Cells.Item(3, 15).Value = 35.737266
$ws.Cells.Item(3, 16).Value = 59.701065
$ws.Cells.Item(3, 17).Value = -0.2
$ws.Cells.Item(3, 18).Value = 0.2
$ws.Cells.Item(3, 19).Value = 0.427646
$ws.Cells.Item(3, 20).Value = 0.179644
$ws.Cells.Item(3, 21).Value = 0.392535
$ws.Cells.Item(3, 22).Value = 0.238485
$ws.Cells.Item(3, 23).Value = 0.761341
$ws.Cells.Item(3, 24).Value = 0.406673
$ws.Cells.Item(3, 25).Value = 0.593153
$ws.Cells.Item(3, 26).Value = 0.581479
$ws.Cells.Item(3, 27).Value = 0.418346
$ws.Cells.Item(3, 28).Value = 0.732884
$ws.Cells.Item(3, 29).Value = 0.266941
$ws.Cells.Item(3, 30).Value = 0.845286
$ws.Cells.Item(3, 31).Value = 0.154539
$ws.Cells.Item(3, 32).Value = 0.742102
$ws.Cells.Item(3, 33).Value = 0.257898
$ws.Cells.Item(3, 34).Value = 0.494021
$ws.Cells.Item(3, 35).Value = 0.505979
$ws.Cells.Item(3, 36).Value = 0.722187
$ws.Cells.Item(3, 37).Value = 0.277813
$ws.Cells.Item(3, 38).Value = 0.468153
$ws.Cells.Item(3, 39).Value = 0.531847
$ws.Cells.Item(3, 40).Value = 0.765397
$ws.Cells.Item(3, 41).Value = 0.736681

# Row 4
$ws.Cells.Item(4, 1).Value = 1369
$ws.Cells.Item(4, 2).Value = "2025-12-01T19:30:00"
$ws.Cells.Item(4, 3).Value = "ЦСКА"
$ws.Cells.Item(4, 4).Value = "Салават Юлаев"
$ws.Cells.Item(4, 5).Value = 897834
$ws.Cells.Item(4, 6).Value = "https://text.khl.ru/text/897834.html"
$ws.Cells.Item(4, 7).Value = 3.561942
$ws.Cells.Item(4, 8).Value = 1.027778
$ws.Cells.Item(4, 9).Value = 1.166667
$ws.Cells.Item(4, 10).Value = 2.157709
$ws.Cells.Item(4, 11).Value = 2.859826
$ws.Cells.Item(4, 12).Value = 1.097222
$ws.Cells.Item(4, 13).Value = 4.58972
$ws.Cells.Item(4, 14).Value = 29.27658
$ws.Cells.Item(4, 15).Value = 21.790523
$ws.Cells.Item(4, 16).Value = 51.067103
$ws.Cells.Item(4, 17).Value = 0.07492600000000001
$ws.Cells.Item(4, 18).Value = -0.2
$ws.Cells.Item(4, 19).Value = 0.736415
$ws.Cells.Item(4, 20).Value = 0.146246
$ws.Cells.Item(4, 21).Value = 0.117143
$ws.Cells.Item(4, 22).Value = 0.441906
$ws.Cells.Item(4, 23).Value = 0.557898
$ws.Cells.Item(4, 24).Value = 0.637228
$ws.Cells.Item(4, 25).Value = 0.362576
$ws.Cells.Item(4, 26).Value = 0.791807
$ws.Cells.Item(4, 27).Value = 0.207997
$ws.Cells.Item(4, 28).Value = 0.893753
$ws.Cells.Item(4, 29).Value = 0.106051
$ws.Cells.Item(4, 30).Value = 0.951383
$ws.Cells.Item(4, 31).Value = 0.048421
$ws.Cells.Item(4, 32).Value = 0.778914
$ws.Cells.Item(4, 33).Value = 0.221086
$ws.Cells.Item(4, 34).Value = 0.5446839999999999
$ws.Cells.Item(4, 35).Value = 0.455316
$ws.Cells.Item(4, 36).Value = 0.299953
$ws.Cells.Item(4, 37).Value = 0.700047
$ws.Cells.Item(4, 38).Value = 0.099025
$ws.Cells.Item(4, 39).Value = 0.900975
$ws.Cells.Item(4, 40).Value = 0.959053
$ws.Cells.Item(4, 41).Value = 0.462498

# Row 5
$ws.Cells.Item(5, 1).Value = 1369
$ws.Cells.Item(5, 2).Value = "2025-12-01T19:30:00"
$ws.Cells.Item(5, 3).Value = "Динамо М"
$ws.Cells.Item(5, 4).Value = "Торпедо"
$ws.Cells.Item(5, 5).Value = 897837
$ws.Cells.Item(5, 6).Value = "https://text.khl.ru/text/897837.html"
$ws.Cells.Item(5, 7).Value = 1.460695
$ws.Cells.Item(5, 8).Value = 2.942447
$ws.Cells.Item(5, 9).Value = 1.503297
$ws.Cells.Item(5, 10).Value = 1.325609
$ws.Cells.Item(5, 11).Value = 1.393152
$ws.Cells.Item(5, 12).Value = 2.222872
$ws.Cells.Item(5, 13).Value = 4.403142
$ws.Cells.Item(5, 14).Value = 24.776157
$ws.Cells.Item(5, 15).Value = 31.069775
$ws.Cells.Item(5, 16).Value = 55.845933
$ws.Cells.Item(5, 17).Value = -0.180896
$ws.Cells.Item(5, 18).Value = 0.02035
$ws.Cells.Item(5, 19).Value = 0.235408
$ws.Cells.Item(5, 20).Value = 0.201686
$ws.Cells.Item(5, 21).Value = 0.5628840000000001
$ws.Cells.Item(5, 22).Value = 0.511816
$ws.Cells.Item(5, 23).Value = 0.488162
$ws.Cells.Item(5, 24).Value = 0.703372
$ws.Cells.Item(5, 25).Value = 0.296606
$ws.Cells.Item(5, 26).Value = 0.841905
$ws.Cells.Item(5, 27).Value = 0.158073
$ws.Cells.Item(5, 28).Value = 0.925396
$ws.Cells.Item(5, 29).Value = 0.074582
$ws.Cells.Item(5, 30).Value = 0.968525
$ws.Cells.Item(5, 31).Value = 0.031453
$ws.Cells.Item(5, 32).Value = 0.405801
$ws.Cells.Item(5, 33).Value = 0.594199
$ws.Cells.Item(5, 34).Value = 0.16485
$ws.Cells.Item(5, 35).Value = 0.8351499999999999
$ws.Cells.Item(5, 36).Value = 0.650971
$ws.Cells.Item(5, 37).Value = 0.349029
$ws.Cells.Item(5, 38).Value = 0.383413
$ws.Cells.Item(5, 39).Value = 0.616587
$ws.Cells.Item(5, 40).Value = 0.651632
$ws.Cells.Item(5, 41).Value = 0.899029
